# Regenerate save_data: update column G (header "K") values for rows 2-31
# (Strike# column was regenerated to hold K values; row 28 is unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    29 = 2
    30 = 0
    31 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
